$wb = $excel.ActiveWorkbook

# Trade #6 closed at 2026-02-16 22:57:50 - base_strategy UP +0.000%
# New trade row appended as row 7 on both the "All Trades" ledger and the
# per-strategy "base_strategy" sheet (they mirror the same trade log).

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 7

    # Keep the date/time columns as literal text (matching the existing
    # rows above) instead of letting them auto-convert to date/time serials.
    $ws.Range("B$row`:C$row").NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value  = 6                      # Trade #
    $ws.Cells.Item($row, 2).Value  = "2026-02-16"            # Date
    $ws.Cells.Item($row, 3).Value  = "22:57:50"              # Time
    $ws.Cells.Item($row, 4).Value  = "base_strategy"         # Strategy
    $ws.Cells.Item($row, 5).Value  = "UP"                    # Side
    $ws.Cells.Item($row, 6).Value  = 0.5                     # Entry Price
    # Column G (Exit Price) left blank
    $ws.Cells.Item($row, 8).Value  = "OPEN"                  # Status
    $ws.Cells.Item($row, 9).Value  = 0                       # P&L %
    $ws.Cells.Item($row, 10).Value = 0                       # P&L $
    $ws.Cells.Item($row, 11).Value = 100                     # Capital After
    $ws.Cells.Item($row, 12).Value = 0                       # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                       # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                     # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    # Column P (Exit Reason) left blank
    $ws.Cells.Item($row, 17).Value = 0                       # Duration (min)
}
